$d = $word.ActiveDocument

$map = @{
    "TC-010" = "TC-006"
    "TC-011" = "TC-007"
    "TC-012" = "TC-008"
    "TC-013" = "TC-009"
    "TC-014" = "TC-010"
    "TC-015" = "TC-011"
    "TC-016" = "TC-012"
    "TC-017" = "TC-013"
    "TC-018" = "TC-014"
    "TC-019" = "TC-015"
    "TC-020" = "TC-016"
    "TC-021" = "TC-017"
    "TC-022" = "TC-18"
    "TC-023" = "TC-019"
    "TC-024" = "TC-020"
    "TC-025" = "TC-021"
    "TC-026" = "TC-022"
    "TC-027" = "TC-023"
    "TC-028" = "TC-024"
    "TC-029" = "TC-025"
    "TC-030" = "TC-026"
    "TC-031" = "TC-027"
    "TC-032" = "TC-028"
    "TC-033" = "TC-029"
    "TC-034" = "TC-030"
    "TC-035" = "TC-031"
    "TC-036" = "TC-032"
    "TC-037" = "TC-033"
    "TC-038" = "TC-034"
    "TC-039" = "TC-035"
    "TC-040" = "TC-036"
    "TC-041" = "TC-037"
    "TC-042" = "TC-038"
    "TC-043" = "TC-039"
    "TC-044" = "TC-040"
    "TC-045" = "TC-041"
    "TC-046" = "TC-042"
    "TC-047" = "TC-043"
    "TC-048" = "TC-044"
    "TC-049" = "TC-045"
    "TC-050" = "TC-046"
    "TC-051" = "TC-047"
    "TC-052" = "TC-048"
    "TC-053" = "TC-049"
    "TC-054" = "TC-050"
    "TC-055" = "TC-051"
}

for ($i = 2; $i -le 11; $i++) {
    $t = $d.Tables.Item($i)
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        $cell = $t.Cell($r, 1)
        $cellRange = $cell.Range
        $txt = $cellRange.Text
        $trimmed = $txt.TrimEnd([char]7, [char]13)
        if ($map.ContainsKey($trimmed)) {
            $cellRange.Text = $map[$trimmed]
        }
    }
}
